# Update trading results - Thu Dec 11 01:46:21 UTC 2025
# Appends 4 new log rows (130-133) to the trading log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 130; A = "2025-12-11T01:46:17.236707"; B = "TRADING_ATTEMPT"; C = "SOL"; D = "UNKNOWN"; E = 133.5144301747269;  K = "ATTEMPT"; L = "Attempting trade 1/2" },
    @{ Row = 131; A = "2025-12-11T01:46:18.963429"; B = "POSITION_FAILED"; C = "SOL"; D = "UNKNOWN"; E = $null;              K = "FAILED";  L = "Trade execution failed for trade 1" },
    @{ Row = 132; A = "2025-12-11T01:46:19.008333"; B = "TRADING_ATTEMPT"; C = "ETH"; D = "UNKNOWN"; E = 3267.627119741301;  K = "ATTEMPT"; L = "Attempting trade 2/2" },
    @{ Row = 133; A = "2025-12-11T01:46:20.481552"; B = "POSITION_FAILED"; C = "ETH"; D = "UNKNOWN"; E = $null;              K = "FAILED";  L = "Trade execution failed for trade 2" }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row

    $ws.Cells.Item($rowIndex, 1).Value = $r.A
    $ws.Cells.Item($rowIndex, 2).Value = $r.B
    $ws.Cells.Item($rowIndex, 3).Value = $r.C
    $ws.Cells.Item($rowIndex, 4).Value = $r.D

    if ($null -ne $r.E) {
        $ws.Cells.Item($rowIndex, 5).Value = $r.E
    } else {
        $ws.Cells.Item($rowIndex, 5).Value = ""
    }

    $ws.Cells.Item($rowIndex, 6).Value = ""
    $ws.Cells.Item($rowIndex, 7).Value = ""
    $ws.Cells.Item($rowIndex, 8).Value = ""
    $ws.Cells.Item($rowIndex, 9).Value = ""
    $ws.Cells.Item($rowIndex, 10).Value = ""

    $ws.Cells.Item($rowIndex, 11).Value = $r.K
    $ws.Cells.Item($rowIndex, 12).Value = $r.L
}
